# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The underlying worker/period database behind this "Estado de Cuenta"
# report changed:
#   - Two workers (OMAR JOSE MADRID YEPEZ / 1044910682 and
#     MONTSERRAT GUDELIA ROMEU BOSSIO / 45566027) were removed.
#   - The three remaining workers (DOMINGO GUERRA GUERRERO, DURERLYS MARIA
#     CUADRO ARRIETA, DIOSMEL ENRIQUE SUAREZ GARCIA) now each carry THREE
#     overdue periods (2506, 2507 and the newly added 2508) instead of two
#     (2506, 2507), and the rows are grouped by period instead of by worker.
#   - The summary figures (total overdue value, worker count, period count)
#     were refreshed to match.
#
# This script reproduces that edit against the already-open workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refreshed summary figures above the table
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 607968   # VALOR MORA total (was 668512)
$ws.Range("C13").Value = 3        # Cant. Trabajadores (was 5)
$ws.Range("F13").Value = 3        # Cant. Periodos (was 2)

# ---------------------------------------------------------------------
# 2) Rewrite the worker/period detail table (rows 16-24).
#    Tipo Doc | N Doc | Nombre | Periodo Mora | Valor Mora | Salario Basico
#
#    The table used to hold 10 data rows (16-25); it now only needs 9
#    (16-24). Remove one row first so everything shifts up by one and
#    the closing (bottom-border) row formatting that lived on the old
#    last row (25) ends up on the new last row (24) - exactly as when
#    the report is regenerated with one fewer row - then fill in the
#    refreshed data top to bottom.
# ---------------------------------------------------------------------
$ws.Rows(16).Delete()

$rows = @(
    @("CC", "73129331",   "DOMINGO GUERRA GUERRERO",       "2506", 67600, 1690000),
    @("CC", "45563706",   "DURERLYS MARIA CUADRO ARRIETA", "2506", 78116, 1952900),
    @("CC", "1235047647", "DIOSMEL ENRIQUE SUAREZ GARCIA", "2506", 56940, 1423500),
    @("CC", "73129331",   "DOMINGO GUERRA GUERRERO",       "2507", 67600, 1690000),
    @("CC", "45563706",   "DURERLYS MARIA CUADRO ARRIETA", "2507", 78116, 1952900),
    @("CC", "1235047647", "DIOSMEL ENRIQUE SUAREZ GARCIA", "2507", 56940, 1423500),
    @("CC", "73129331",   "DOMINGO GUERRA GUERRERO",       "2508", 67600, 1690000),
    @("CC", "45563706",   "DURERLYS MARIA CUADRO ARRIETA", "2508", 78116, 1952900),
    @("CC", "1235047647", "DIOSMEL ENRIQUE SUAREZ GARCIA", "2508", 56940, 1423500)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("B$r").Value = $data[0]
    $ws.Range("C$r").Value = $data[1]
    $ws.Range("D$r").Value = $data[2]
    $ws.Range("E$r").Value = $data[3]
    $ws.Range("F$r").Value = $data[4]
    $ws.Range("G$r").Value = $data[5]
}

# Column D ("Nombre Trabajador") was sized to fit the longest name that
# is no longer present; let it re-fit the remaining content.
$ws.Columns("D").AutoFit()
